$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 4 entirely (it's being removed)
$ws.Range("A4:F4").ClearContents()

# Header row (row 1)
$ws.Range("A1").Value = "SAN NICOLAS DISANG S.R.L"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"
$ws.Range("D1").Value = "Unnamed: 3"
$ws.Range("E1").Value = "Unnamed: 4"
$ws.Range("F1").Value = "Unnamed: 5"
$ws.Range("G1").Value = "Unnamed: 6"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Row 2
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "POXILINA"
$ws.Range("C2").Value = "10MIN."
$ws.Range("D2").Value = "70G"
$ws.Range("E2").Value = 1431.41
$ws.Range("F2").Value = 889.46
$ws.Range("G2").Value = 1076.25

# Row 3
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "POXILINA"
$ws.Range("C3").Value = "10MIN."
$ws.Range("D3").Value = "250G"
$ws.Range("E3").Value = 3138.62
$ws.Range("F3").Value = 1950.3
$ws.Range("G3").Value = 2359.86
